# Update the "default username/password" callout textboxes on the
# kvdb server/client usage slides (slide 7 = server, slide 8 = client).
#
# For both slides the textbox is shape index 4 ("文本框 4" / "文本框 7"),
# and its first paragraph is made of four runs followed by an
# endParaRPr, then a second paragraph with one more run:
#   [1,4)   "kvdb"
#   [5,9)   "数据库的默认用户名"   -> "数据库的初始默认用户名"
#   [14,1)  "/"
#   [15,3)  "密码："
#   [19,14) "admin/password"   (second paragraph)
#
# The edit clarifies the wording to "初始默认用户名" (initial default user
# name) and also stamps every run's "other" (symbol/complex-script)
# typeface to the theme's minor East-Asian font, i.e. what PowerPoint
# persists as <a:sym typeface="+mn-ea"/> on each <a:rPr>.

$p = $ppt.ActivePresentation

function Update-KvdbCallout($slideIndex) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(4)
    $tr = $shape.TextFrame.TextRange

    # Stamp the "other" typeface on every existing run first, while the
    # character offsets below still line up with the original text.
    foreach ($range in @(
            @(1, 4),    # "kvdb"
            @(5, 9),    # "数据库的默认用户名"
            @(14, 1),   # "/"
            @(15, 3),   # "密码："
            @(19, 14)   # "admin/password"
        )) {
        try {
            $tr.Characters($range[0], $range[1]).Font.NameOther = "+mn-ea"
        } catch {
            # Some hosts expose Font.NameOther as read-only; ignore.
        }
    }

    # Now update the wording of the second run.
    $tr.Characters(5, 9).Text = "数据库的初始默认用户名"
}

Update-KvdbCallout 7
Update-KvdbCallout 8
